$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.144.80'
$ws.Range('E2').Value = '  -4.82%  '
$ws.Range('D3').Value = '3.252.24'
$ws.Range('E3').Value = '  -6.10%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.30'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.91%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.584'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.82%  '
$ws.Range('D9').Value = '3.242.39'
$ws.Range('E9').Value = '  -5.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.183'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.581'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '46.98'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -8.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000262'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -7.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '629.32'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.48'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.96%  '
$ws.Range('D16').Value = '3.792.73'
$ws.Range('E16').Value = '  -6.07%  '
$ws.Range('D17').Value = '65.229.24'
$ws.Range('E17').Value = '  -4.77%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.116'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.32%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.67'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.76%  '
$ws.Range('D20').Value = '3.277.56'
$ws.Range('E20').Value = '  -6.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.26'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -8.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.897'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.56'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '105.62'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.95'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.38'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.63'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.04'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.92'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.26'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.10%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '554.74'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +11.72%  '
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.99'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.104'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '57.19'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.91'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +38.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').Value = '3.592.21'
$ws.Range('E39').Value = '  -2.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.45'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.69'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.21%  '
$ws.Range('D42').Value = '0.0₃0700'
$ws.Range('E42').Value = '  -9.91%  '
$ws.Range('E43').Value = '  -4.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.338'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '31.71'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.93%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.24'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.15%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0411'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.55%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.128'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.81%  '
$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.58'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.24'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.61%  '
